$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1001
$ws.Range("B1").Value = "John Wick"
$ws.Range("C1").Value = 89

$ws.Range("A2").Value = 1002
$ws.Range("B2").Value = "James Bond"
$ws.Range("C2").Value = 82
